# Mixorcerers TODO workbook update:
#  - Backlog!B4 ("Write Local Server") is re-highlighted to match the
#    "Client wireframe to test server" task's green fill.
#  - Backlog!B6's finished task text ("Create Local Server Scene") is
#    cleared, the cell keeping a pale "done" highlight.
#  - A new backlog item "Low  processor mode in menus" (priority 5) is
#    added as row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# 1) Re-color B4 to match B5's green highlight (reuse existing format
#    exactly via copy/paste so no new style entries are created).
$ws.Range("B5").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# 2) Clear the now-finished "Create Local Server Scene" task text from
#    B6, leaving a pale highlight behind to mark it done.
$ws.Range("B6").Interior.Color = 16777215
$ws.Range("B6").Interior.PatternColor = 13434879
$ws.Range("B6").ClearContents()

# 3) Add the new backlog entry as row 11.
$ws.Range("A11").Value = "Low  processor mode in menus"
$ws.Range("C11").Value = 5

$ws.Range("C11").Select()
